# "Generate Report for Archive" - refresh the localization status report:
# flip the in-flight items from "Ready for handoff" to "In Translation"
# on the Overview rollup sheet and on each per-locale detail sheet, then
# re-fit the affected columns to the new text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# The status text got shorter, so the columns that hold it are re-fit to
# the new (narrower) content.
$refitColumnWidth = 12.5

# Overview sheet: "zh-cn" and "de-de" status columns (E and F)
$wsOverview = $wb.Worksheets.Item("Overview")
$usedRows = $wsOverview.UsedRange.Rows.Count
for ($r = 2; $r -le $usedRows; $r++) {
    foreach ($col in @("E", "F")) {
        $cell = $wsOverview.Range($col + $r)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value2 = $newStatus
        }
    }
}
$wsOverview.Columns("E:F").AutoFit() | Out-Null
$wsOverview.Columns("E:F").ColumnWidth = $refitColumnWidth

# Per-locale detail sheets: "Status" column (C)
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $usedRows = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $usedRows; $r++) {
        $cell = $ws.Range("C" + $r)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value2 = $newStatus
        }
    }
    $ws.Columns("C:C").AutoFit() | Out-Null
    $ws.Columns("C:C").ColumnWidth = $refitColumnWidth
}
